# Fix 무학 (Moohak) IFRS financial data rows 2-9 (error solve ifrs list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2901
$ws.Range("E2").Value = 814
$ws.Range("F2").Value = 814
$ws.Range("G2").Value = 1082
$ws.Range("H2").Value = 829
$ws.Range("I2").Value = 829
$ws.Range("K2").Value = 5342
$ws.Range("L2").Value = 1192
$ws.Range("M2").Value = 4151
$ws.Range("N2").Value = 4151
$ws.Range("P2").Value = 55
$ws.Range("Q2").Value = 1283
$ws.Range("R2").Value = -1265
$ws.Range("S2").Value = -32
$ws.Range("T2").Value = 410
$ws.Range("U2").Value = 874
$ws.Range("W2").Value = 28.07
$ws.Range("X2").Value = 28.56
$ws.Range("Y2").Value = 21.77
$ws.Range("Z2").Value = 16.99
$ws.Range("AA2").Value = 28.71
$ws.Range("AB2").Value = 7251.55
$ws.Range("AC2").Value = 2908
$ws.Range("AD2").Value = 12.74
$ws.Range("AE2").Value = 14750
$ws.Range("AF2").Value = 2.51
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 28500000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("V2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2958
$ws.Range("E3").Value = 657
$ws.Range("F3").Value = 657
$ws.Range("G3").Value = 433
$ws.Range("H3").Value = 288
$ws.Range("I3").Value = 288
$ws.Range("K3").Value = 5568
$ws.Range("L3").Value = 1209
$ws.Range("M3").Value = 4359
$ws.Range("N3").Value = 4359
$ws.Range("P3").Value = 56
$ws.Range("Q3").Value = 783
$ws.Range("R3").Value = -649
$ws.Range("S3").Value = -82
$ws.Range("T3").Value = 373
$ws.Range("U3").Value = 410
$ws.Range("W3").Value = 22.2
$ws.Range("X3").Value = 9.74
$ws.Range("Y3").Value = 6.77
$ws.Range("Z3").Value = 5.28
$ws.Range("AA3").Value = 27.73
$ws.Range("AB3").Value = 7607.26
$ws.Range("AC3").Value = 1011
$ws.Range("AD3").Value = 37.79
$ws.Range("AE3").Value = 15603
$ws.Range("AF3").Value = 2.45
$ws.Range("AG3").Value = 246
$ws.Range("AH3").Value = 0.65
$ws.Range("AI3").Value = 23.89
$ws.Range("AJ3").Value = 28500000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("V3").ClearContents()

# Row 4
$ws.Range("D4").Value = 2702
$ws.Range("E4").Value = 520
$ws.Range("F4").Value = 520
$ws.Range("G4").Value = 821
$ws.Range("H4").Value = 615
$ws.Range("I4").Value = 615
$ws.Range("K4").Value = 6037
$ws.Range("L4").Value = 1128
$ws.Range("M4").Value = 4910
$ws.Range("N4").Value = 4910
$ws.Range("P4").Value = 57
$ws.Range("Q4").Value = 528
$ws.Range("R4").Value = -257
$ws.Range("S4").Value = -116
$ws.Range("T4").Value = 288
$ws.Range("U4").Value = 240
$ws.Range("W4").Value = 19.23
$ws.Range("X4").Value = 22.76
$ws.Range("Y4").Value = 13.27
$ws.Range("Z4").Value = 10.6
$ws.Range("AA4").Value = 22.97
$ws.Range("AB4").Value = 8457.9
$ws.Range("AC4").Value = 2157
$ws.Range("AD4").Value = 10.59
$ws.Range("AE4").Value = 17542
$ws.Range("AF4").Value = 1.3
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 1.53
$ws.Range("AI4").Value = 15.93
$ws.Range("AJ4").Value = 28500000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()

# Row 5
$ws.Range("D5").Value = 2505
$ws.Range("E5").Value = 287
$ws.Range("F5").Value = 287
$ws.Range("G5").Value = 771
$ws.Range("H5").Value = 515
$ws.Range("I5").Value = 515
$ws.Range("K5").Value = 6701
$ws.Range("L5").Value = 1343
$ws.Range("M5").Value = 5358
$ws.Range("N5").Value = 5358
$ws.Range("P5").Value = 57
$ws.Range("Q5").Value = 196
$ws.Range("R5").Value = -342
$ws.Range("S5").Value = 202
$ws.Range("T5").Value = 214
$ws.Range("U5").Value = -18
$ws.Range("V5").Value = 300
$ws.Range("W5").Value = 11.46
$ws.Range("X5").Value = 20.56
$ws.Range("Y5").Value = 10.03
$ws.Range("Z5").Value = 8.09
$ws.Range("AA5").Value = 25.07
$ws.Range("AB5").Value = 9176.950000000001
$ws.Range("AC5").Value = 1807
$ws.Range("AD5").Value = 10.87
$ws.Range("AE5").Value = 19008
$ws.Range("AF5").Value = 1.03
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 1.78
$ws.Range("AI5").Value = 19.15
$ws.Range("AJ5").Value = 28500000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 1937
$ws.Range("E6").Value = -100
$ws.Range("F6").Value = -100
$ws.Range("G6").Value = -262
$ws.Range("H6").Value = -203
$ws.Range("I6").Value = -203
$ws.Range("K6").Value = 6190
$ws.Range("L6").Value = 1125
$ws.Range("M6").Value = 5064
$ws.Range("N6").Value = 5064
$ws.Range("P6").Value = 57
$ws.Range("Q6").Value = -177
$ws.Range("R6").Value = -144
$ws.Range("S6").Value = 21
$ws.Range("T6").Value = 297
$ws.Range("U6").Value = -473
$ws.Range("V6").Value = 420
$ws.Range("W6").Value = -5.18
$ws.Range("X6").Value = -10.46
$ws.Range("Y6").Value = -3.89
$ws.Range("Z6").Value = -3.14
$ws.Range("AA6").Value = 22.22
$ws.Range("AB6").Value = 8701.16
$ws.Range("AC6").Value = -711
$ws.Range("AD6").Value = -19.13
$ws.Range("AE6").Value = 17967
$ws.Range("AF6").Value = 0.76
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 2.57
$ws.Range("AI6").Value = -48.69
$ws.Range("AJ6").Value = 28500000

# Row 7
$ws.Range("D7").Value = 1860
$ws.Range("E7").Value = -90
$ws.Range("I7").Value = -190
$ws.Range("W7").Value = -4.84
$ws.Range("AC7").Value = -667
$ws.Range("AD7").Value = -12.35
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 1817
$ws.Range("E8").Value = 13
$ws.Range("I8").Value = -111
$ws.Range("W8").Value = 0.72
$ws.Range("AC8").Value = -389
$ws.Range("AD8").Value = -21.13
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 1880
$ws.Range("E9").Value = 20
$ws.Range("I9").Value = -110
$ws.Range("W9").Value = 1.06
$ws.Range("AC9").Value = -386
$ws.Range("AD9").Value = -21.32
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

Write-Host "무학 IFRS data updated."